$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parametros")

# Insert a new column before column D (Bus m stays C, new "id" column becomes D,
# old D..H (r,x,b,Rating,Costo) shift right to E..I)
$ws.Range("D1").EntireColumn.Insert()

# Header for the new column, matching the style of the other header cells
$ws.Range("D1").Value = "id"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("D1").VerticalAlignment = -4160
$ws.Range("D1").Borders.LineStyle = 1

# Data values for the new column
$ws.Range("D2").Value = "NL"
$ws.Range("D3").Value = "NL"

# Update selection to match the final state
$ws.Range("D4").Select()
